$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.144.22"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.67%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.831.19"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.72%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.009"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.95%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.70%  "
$ws.Range("E6").Value = "  +0.66%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4703"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.96%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3686"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.49%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07397"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.22%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8807"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.04%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.41"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.26%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.823.98"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.56%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07310"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.466"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.95%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.75"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.46%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.550"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.53%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.009"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.89%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008769"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.52%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.007"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.48%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.77"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "27.166.65"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.71%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.303"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.98%  "
$ws.Range("E23").Value = "  +1.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.069.30"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.43%  "
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("E26").Value = "  +0.25%  "
$ws.Range("E27").Value = "  +0.48%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.155"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.36%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.267"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.89%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "117.39"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.58%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08922"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7589"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.59%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.170"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.53%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.541"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.12%  "
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("E36").Value = "  +0.64%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.102"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.54%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05331"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.31%  "
$ws.Range("E39").Value = "  -0.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.000"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.66%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.416"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.47%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.302"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.76%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5341"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.80%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1662"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.547"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.40%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4943"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.48%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.55"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.89%  "
$ws.Range("E48").Value = "  +0.64%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.669"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "103.74"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.52%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06316"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.36%  "
